$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 198.21428
$ws.Range("I39").Value = 73.27273
$ws.Range("J39").Value = 656.3333
$ws.Range("K39").Value = 219.81819
$ws.Range("L39").Value = 1968.9999
$ws.Range("M39").Value = 76.18181000000001
$ws.Range("N39").Value = -2560.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4470.773
$ws.Range("I62").Value = 1334.0625
$ws.Range("J62").Value = 12835.333
$ws.Range("K62").Value = 1334.0625
$ws.Range("L62").Value = 12835.333
$ws.Range("M62").Value = -710.0625
$ws.Range("N62").Value = -14083.333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4470.773
$ws.Range("I65").Value = 1334.0625
$ws.Range("J65").Value = 12835.333
$ws.Range("K65").Value = 6670.3125
$ws.Range("L65").Value = 64176.665
$ws.Range("M65").Value = -3550.3125
$ws.Range("N65").Value = -70416.66500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 628.4
$ws.Range("I99").Value = 350.5
$ws.Range("J99").Value = 1740
$ws.Range("K99").Value = 1051.5
$ws.Range("L99").Value = 5220
$ws.Range("M99").Value = 446.5
$ws.Range("N99").Value = -8216

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 63780
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 63780
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 63780
$ws.Range("N123").Value = -73580

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1118.52
$ws.Range("I129").Value = 498.8
$ws.Range("J129").Value = 1273.45
$ws.Range("K129").Value = 1496.4
$ws.Range("L129").Value = 3820.35
$ws.Range("M129").Value = 3503.6
$ws.Range("N129").Value = -13820.35

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 7233.091
$ws.Range("I131").Value = 1583
$ws.Range("J131").Value = 22300
$ws.Range("K131").Value = 4749
$ws.Range("L131").Value = 66900
$ws.Range("M131").Value = 291
$ws.Range("N131").Value = -76980

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4359554
$ws.Range("I138").Value = 1610839
$ws.Range("J138").Value = 5379238.5
$ws.Range("K138").Value = 4832517
$ws.Range("L138").Value = 16137715.5
$ws.Range("M138").Value = -4827377
$ws.Range("N138").Value = -16147995.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1564.4736
$ws.Range("I141").Value = 1164.7059
$ws.Range("J141").Value = 4962.5
$ws.Range("K141").Value = 3494.1177
$ws.Range("L141").Value = 14887.5
$ws.Range("M141").Value = 1685.8823
$ws.Range("N141").Value = -25247.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1935.6207
$ws.Range("I61").Value = 1224.7391
$ws.Range("J61").Value = 4660.6665
$ws.Range("K61").Value = 1224.7391
$ws.Range("L61").Value = 4660.6665
$ws.Range("M61").Value = -1012.7391
$ws.Range("N61").Value = -5084.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1531.5385
$ws.Range("I102").Value = 1446.3636
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1446.3636
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 175.6364000000001
$ws.Range("N102").Value = -5244

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1562.1708
$ws.Range("I122").Value = 1270.1154
$ws.Range("J122").Value = 2068.4
$ws.Range("K122").Value = 3810.3462
$ws.Range("L122").Value = 6205.200000000001
$ws.Range("M122").Value = -1360.3462
$ws.Range("N122").Value = -11105.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1935.6207
$ws.Range("I136").Value = 1224.7391
$ws.Range("J136").Value = 4660.6665
$ws.Range("K136").Value = 3674.2173
$ws.Range("L136").Value = 13981.9995
$ws.Range("M136").Value = -1124.2173
$ws.Range("N136").Value = -19081.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 13164
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 13164
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 13164
$ws.Range("N103").Value = -15508

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3292.9666
$ws.Range("I105").Value = 3035.8635
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3035.8635
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -1288.8635
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1496.4054
$ws.Range("I31").Value = 959.2143
$ws.Range("J31").Value = 3167.6667
$ws.Range("K31").Value = 959.2143
$ws.Range("L31").Value = 3167.6667
$ws.Range("M31").Value = -664.2143
$ws.Range("N31").Value = -3757.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 21459.6
$ws.Range("I33").Value = 21459.6
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 21459.6
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -21080.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1496.4054
$ws.Range("I34").Value = 959.2143
$ws.Range("J34").Value = 3167.6667
$ws.Range("K34").Value = 959.2143
$ws.Range("L34").Value = 3167.6667
$ws.Range("M34").Value = -757.2143
$ws.Range("N34").Value = -3571.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 24041.4
$ws.Range("I62").Value = 36069.168
$ws.Range("J62").Value = 5999.75
$ws.Range("K62").Value = 36069.168
$ws.Range("L62").Value = 5999.75
$ws.Range("M62").Value = -35445.168
$ws.Range("N62").Value = -7247.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 24041.4
$ws.Range("I65").Value = 36069.168
$ws.Range("J65").Value = 5999.75
$ws.Range("K65").Value = 180345.84
$ws.Range("L65").Value = 29998.75
$ws.Range("M65").Value = -177225.84
$ws.Range("N65").Value = -36238.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 20834400
$ws.Range("I99").Value = 20834400
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 20834400
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -20832902
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 939.1
$ws.Range("I105").Value = 932.3333
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 932.3333
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 814.6667
$ws.Range("N105").Value = -4494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2568.375
$ws.Range("I122").Value = 1166.4
$ws.Range("J122").Value = 4905
$ws.Range("K122").Value = 3499.2
$ws.Range("L122").Value = 14715
$ws.Range("M122").Value = -1049.2
$ws.Range("N122").Value = -19615

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 20834400
$ws.Range("I126").Value = 20834400
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 62503200
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -62500730
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 13040
$ws.Range("I116").Value = 13040
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 39120
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -35678
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1031
$ws.Range("I122").Value = 502
$ws.Range("J122").Value = 1242.6
$ws.Range("K122").Value = 4518
$ws.Range("L122").Value = 11183.4
$ws.Range("M122").Value = -2068
$ws.Range("N122").Value = -16083.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1282.3572
$ws.Range("I129").Value = 553.3333
$ws.Range("J129").Value = 1829.125
$ws.Range("K129").Value = 1659.9999
$ws.Range("L129").Value = 5487.375
$ws.Range("M129").Value = 3340.0001
$ws.Range("N129").Value = -15487.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1328.6753
$ws.Range("I131").Value = 354.63635
$ws.Range("J131").Value = 1491.0151
$ws.Range("K131").Value = 1063.90905
$ws.Range("L131").Value = 4473.0453
$ws.Range("M131").Value = 3976.09095
$ws.Range("N131").Value = -14553.0453

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1435.8
$ws.Range("I132").Value = 850
$ws.Range("J132").Value = 1500.8889
$ws.Range("K132").Value = 7650
$ws.Range("L132").Value = 13508.0001
$ws.Range("M132").Value = -5120
$ws.Range("N132").Value = -18568.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4508.6855
$ws.Range("I140").Value = 5915.95
$ws.Range("J140").Value = 2632.3333
$ws.Range("K140").Value = 17747.85
$ws.Range("L140").Value = 7896.999899999999
$ws.Range("M140").Value = -12567.85
$ws.Range("N140").Value = -18256.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1579.2122
$ws.Range("I102").Value = 1454.3478
$ws.Range("J102").Value = 1866.4
$ws.Range("K102").Value = 1454.3478
$ws.Range("L102").Value = 1866.4
$ws.Range("M102").Value = 167.6522
$ws.Range("N102").Value = -5110.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3264.7058
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 3535.7144
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 3535.7144
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -3759.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3576.0476
$ws.Range("I122").Value = 1774.75
$ws.Range("J122").Value = 3999.8823
$ws.Range("K122").Value = 5324.25
$ws.Range("L122").Value = 11999.6469
$ws.Range("M122").Value = -2874.25
$ws.Range("N122").Value = -16899.6469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3264.7058
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3535.7144
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 10607.1432
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -15547.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3500.0588
$ws.Range("I132").Value = 2518.1177
$ws.Range("J132").Value = 5463.9414
$ws.Range("K132").Value = 7554.353099999999
$ws.Range("L132").Value = 16391.8242
$ws.Range("M132").Value = -5024.353099999999
$ws.Range("N132").Value = -21451.8242

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 61463.53
$ws.Range("I126").Value = 73992.86
$ws.Range("J126").Value = 2993.3333
$ws.Range("K126").Value = 221978.58
$ws.Range("L126").Value = 8979.999899999999
$ws.Range("M126").Value = -219508.58
$ws.Range("N126").Value = -13919.9999
